$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new name text (literal, including surrounding quotes), new C value.
# The name text begins with a literal apostrophe; escape it as '' so Excel stores
# the leading quote as real content rather than treating it as a quote-prefix marker.
$data = @(
  @(2,  "''Bacteroides_cellulosilyticus_DSM_14838.mat'", 0),
  @(3,  "''Bacteroides_coprocola_M16_DSM_17136.mat'", 0),
  @(4,  "''Bacteroides_fluxus_YIT_12057.mat'", 0),
  @(5,  "''Bacteroides_oleiciplenus_YIT_12058.mat'", 0),
  @(6,  "''Bacteroides_ovatus_ATCC_8483.mat'", 0),
  @(7,  "''Bacteroides_plebeius_M12_DSM_17135.mat'", 0),
  @(8,  "''Bacteroides_salyersiae_WAL_10018.mat'", 0),
  @(9,  "''Bacteroides_stercoris_ATCC_43183.mat'", 0),
  @(10, "''Bacteroides_thetaiotaomicron_VPI_5482.mat'", 0.016),
  @(11, "''Bacteroides_uniformis_ATCC_8492.mat'", 0),
  @(12, "''Bacteroides_vulgatus_ATCC_8482.mat'", 0.01),
  @(13, "''Bifidobacterium_animalis_lactis_AD011.mat'", 0),
  @(14, "''Enterococcus_faecalis_OG1RF_ATCC_47077.mat'", 0),
  @(15, "''Flavonifractor_plautii_ATCC_29863.mat'", 0),
  @(16, "''Lactobacillus_plantarum_JDM1.mat'", 0.014),
  @(17, "''Odoribacter_laneus_YIT_12061.mat'", 0.145),
  @(18, "''Parabacteroides_distasonis_ATCC_8503.mat'", -0),
  @(19, "''Parabacteroides_johnsonii_DSM_18315.mat'", 0.8149999999999999)
)

foreach ($row in $data) {
    $r = $row[0]
    $name = $row[1]
    $val = $row[2]

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $name
    # Setting the value with a leading apostrophe marks the cell as
    # "quote prefixed" (adds a style). Reset style to Normal so the cell
    # keeps no explicit style attribute, same as the original file.
    $cellB.Style = "Normal"

    $cellC = $ws.Cells.Item($r, 3)
    $cellC.Value = $val
}
